$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: rows 1-8 become 0
$ws.Range("A1:A8").Value = 0

# Column I: rows 1-8 become 1
$ws.Range("I1:I8").Value = 1

# H2 gets the 0.2 value that used to live in I2
$ws.Range("H2").Value = 0.2

# Update the active selection to D15
$ws.Range("D15").Select()
